$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A37").Copy()
$ws.Range("A38").PasteSpecial(-4122)
$ws.Range("A38").Value = 45942

$ws.Range("B38").Value = "21,6987"
$ws.Range("C38").Value = "15,6648"
$ws.Range("D38").Value = "15,4517"
$ws.Range("E38").Value = "15,4517"
